# Auto-generated edit script: update cryptos list values (price + volume) per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.260.11"
$ws.Range("D3").Value = "3.561.49"
$ws.Range("E3").Value = "  +1.80%  "
$ws.Range("D5").Value = "'619.21"
$ws.Range("E5").Value = "  +2.88%  "
$ws.Range("E6").Value = "  +3.66%  "
$ws.Range("D7").Value = "3.560.15"
$ws.Range("E7").Value = "  +1.78%  "
$ws.Range("E9").Value = "  +2.31%  "
$ws.Range("E10").Value = "  +5.74%  "
$ws.Range("E11").Value = "  +7.67%  "
$ws.Range("E12").Value = "  +3.81%  "
$ws.Range("E13").Value = "  +2.66%  "
$ws.Range("D14").Value = "'33.24"
$ws.Range("E14").Value = "  +5.46%  "
$ws.Range("D15").Value = "4.163.34"
$ws.Range("E15").Value = "  +1.85%  "
$ws.Range("D16").Value = "3.560.77"
$ws.Range("E16").Value = "  +1.93%  "
$ws.Range("D17").Value = "68.252.67"
$ws.Range("E17").Value = "  +1.48%  "
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("E19").Value = "  +5.71%  "
$ws.Range("E20").Value = "  +6.76%  "
$ws.Range("D21").Value = "'10.06"
$ws.Range("E21").Value = "  +11.81%  "
$ws.Range("D22").Value = "'454.21"
$ws.Range("E22").Value = "  +1.82%  "
$ws.Range("E23").Value = "  +4.29%  "
$ws.Range("D24").Value = "'78.46"
$ws.Range("D25").Value = "'0.0000132"
$ws.Range("E25").Value = "  +2.86%  "
$ws.Range("D26").Value = "3.704.87"
$ws.Range("E26").Value = "  +1.90%  "
$ws.Range("E27").Value = "  -0.11%  "
$ws.Range("D28").Value = "'9.27"
$ws.Range("E28").Value = "  +13.23%  "
$ws.Range("D29").Value = "'10.51"
$ws.Range("E29").Value = "  +4.03%  "
$ws.Range("E30").Value = "  +11.37%  "
$ws.Range("E31").Value = "  +3.67%  "
$ws.Range("E32").Value = "  +4.09%  "
$ws.Range("E34").Value = "  +5.44%  "
$ws.Range("D35").Value = "'26.12"
$ws.Range("E35").Value = "  +1.73%  "
$ws.Range("E36").Value = "  +4.78%  "
$ws.Range("D37").Value = "3.555.61"
$ws.Range("E37").Value = "  +1.97%  "
$ws.Range("E38").Value = "  +3.42%  "
$ws.Range("D39").Value = "'2.38"
$ws.Range("E39").Value = "  +8.65%  "
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("D41").Value = "'181.46"
$ws.Range("E41").Value = "  +4.05%  "
$ws.Range("D42").Value = "'0.0918"
$ws.Range("E42").Value = "  +5.07%  "
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("D44").Value = "'5.65"
$ws.Range("E44").Value = "  +4.88%  "
$ws.Range("D45").Value = "'31.19"
$ws.Range("E45").Value = "  +14.57%  "
$ws.Range("E46").Value = "  +2.15%  "
$ws.Range("D47").Value = "'46.24"
$ws.Range("E47").Value = "  +1.74%  "
$ws.Range("E48").Value = "  +5.79%  "
$ws.Range("D49").Value = "'2.67"
$ws.Range("E49").Value = "  +4.57%  "
$ws.Range("E50").Value = "  +3.43%  "
$ws.Range("E51").Value = "  +7.67%  "
